$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("g10.1")

$ws.Range("B2").Value = -6.943942466780317
$ws.Range("C2").Value = 7.985561900211491
$ws.Range("D2").Value = 3.558416437559409

$ws.Range("B3").Value = -3.033528379010808
$ws.Range("C3").Value = 4.334264824683731
$ws.Range("D3").Value = -5.228869536197389

$ws.Range("B4").Value = -0.2327678074659945
$ws.Range("C4").Value = 4.101953574790462
$ws.Range("D4").Value = 1.61910220795054

$ws.Range("B5").Value = 1.249304374546467
$ws.Range("C5").Value = -1.204086634669399
$ws.Range("D5").Value = 8.070905286572305

$ws.Range("B6").Value = -5.140647810975807
$ws.Range("C6").Value = -3.65206940235947
$ws.Range("D6").Value = 0.1819403087696347

$ws.Range("B7").Value = -3.399975836171742
$ws.Range("C7").Value = 0.5772764519954787
$ws.Range("D7").Value = 0.7720863564637304

$ws.Range("B8").Value = -3.189489044441729
$ws.Range("C8").Value = -1.088355958105625
$ws.Range("D8").Value = -1.90857259036582

$ws.Range("B9").Value = 2.962637469059248
$ws.Range("C9").Value = 1.03338185358528
$ws.Range("D9").Value = 9.619155300664595

$ws.Range("B10").Value = -13.97437319254418
$ws.Range("C10").Value = -3.916628966280844
$ws.Range("D10").Value = -12.16490384706811

$ws.Range("B11").Value = -11.06390641573395
$ws.Range("C11").Value = 15.37912276036422
$ws.Range("D11").Value = -14.12152573771694

$ws.Range("B12").Value = -4.756691349375375
$ws.Range("C12").Value = 14.54125736551517
$ws.Range("D12").Value = -12.94678802350021

$ws.Range("B13").Value = -6.006456415604633
$ws.Range("C13").Value = 7.559924425551756
$ws.Range("D13").Value = -7.503098588847368

$wb.Save()
